# Apply the commit's changes to the workbook:
#  1. On the "Sheet4" tab (matrix size calc data), change the run-length
#     input (D2) and the batch-size inputs (J8/K8) for the row-8 scenario.
#     All the dependent formula cells recalculate automatically.
#  2. Insert a new worksheet ("Sheet6") before "Sheet3", containing two
#     small helper calculations, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Update inputs on Sheet4 -------------------------------------------
$dataSheet = $wb.Worksheets.Item("Sheet4")
$dataSheet.Range("D2").Value = 33
$dataSheet.Range("J8").Value = 1
$dataSheet.Range("K8").Value = 231

# --- 2. Add the new "Sheet6" worksheet, placed right before "Sheet3" ------
$sheet3 = $wb.Worksheets.Item("Sheet3")
$newSheet = $wb.Worksheets.Add($sheet3)

$newSheet.Range("A1").Formula = "=24*3*60"
$newSheet.Range("A2").Formula = "=1100/50"

# Make the newly inserted sheet the active tab, matching the workbook's
# saved "activeTab" view state.
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 100
